$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Row 2 (acc5b70d...) and Row 3 (cccb50be...) move from "Ready for
#     handoff" to "Handed back: in sync with en-US".
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     get populated with hyperlinks mirroring the existing handoff links
#     (A -> E, C -> F).
#   - The "Latest Handback DateTime" (G) column is stamped with the
#     handback time.
# ---------------------------------------------------------------------------

$statusHandedBack = "Handed back: in sync with en-US"

$locales = @(
    @{
        SheetName   = "zh-cn"
        MdTarget    = "https://github.com/OpenLocalizationTest/oltest/blob/3c977fbf9c7ef197389ab3ffb919a8852e631f65/e2e"
        Row2XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efa8e208f7f0e60ef195d83cccb994ee9e757b5f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.zh-cn.xlf"
        Row2XlfName = "acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.zh-cn.xlf"
        Row3XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efa8e208f7f0e60ef195d83cccb994ee9e757b5f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.zh-cn.xlf"
        Row3XlfName = "cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.zh-cn.xlf"
        HandbackTime = "2016-03-08 23:27:34"
    },
    @{
        SheetName   = "de-de"
        MdTarget    = "https://github.com/OpenLocalizationTest/oltest/blob/3c977fbf9c7ef197389ab3ffb919a8852e631f65/e2e"
        Row2XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1891f52aca8f55a104b429ee0bed2dcec3b8cbb5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.de-de.xlf"
        Row2XlfName = "acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.de-de.xlf"
        Row3XlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1891f52aca8f55a104b429ee0bed2dcec3b8cbb5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.de-de.xlf"
        Row3XlfName = "cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.de-de.xlf"
        HandbackTime = "2016-03-08 23:28:02"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)

    # ---- Row 2 : acc5b70d-3d70-4ed8-81fb-134c95d14776.md -----------------
    $ws.Range("B2").Value = $statusHandedBack

    $ws.Hyperlinks.Add(
        $ws.Range("E2"),
        ($locale.MdTarget + "/acc5b70d-3d70-4ed8-81fb-134c95d14776.md"),
        "",
        "",
        "acc5b70d-3d70-4ed8-81fb-134c95d14776.md")
    $ws.Range("E2").Style = "HyperLink"

    $ws.Hyperlinks.Add(
        $ws.Range("F2"),
        $locale.Row2XlfUrl,
        "",
        "",
        $locale.Row2XlfName)
    $ws.Range("F2").Style = "HyperLink"

    $ws.Range("G2").Value = $locale.HandbackTime
    $ws.Range("H2").Value = "Include"

    # ---- Row 3 : cccb50be-7698-471e-8f8f-50e3f32d2e44.md -----------------
    $ws.Range("B3").Value = $statusHandedBack

    $ws.Hyperlinks.Add(
        $ws.Range("E3"),
        ($locale.MdTarget + "/cccb50be-7698-471e-8f8f-50e3f32d2e44.md"),
        "",
        "",
        "cccb50be-7698-471e-8f8f-50e3f32d2e44.md")
    $ws.Range("E3").Style = "HyperLink"

    $ws.Hyperlinks.Add(
        $ws.Range("F3"),
        $locale.Row3XlfUrl,
        "",
        "",
        $locale.Row3XlfName)
    $ws.Range("F3").Style = "HyperLink"

    $ws.Range("G3").Value = $locale.HandbackTime
    $ws.Range("H3").Value = "Include"
}
